# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the regenerated data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5668
$ws1.Range("F4").Value = 643
$ws1.Range("F6").Value = 847
$ws1.Range("F7").Value = 61
$ws1.Range("F8").Value = 378
$ws1.Range("F9").Value = 7

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 51

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5668
$ws4.Range("F4").Value = 643
$ws4.Range("F6").Value = 847
$ws4.Range("F7").Value = 61
$ws4.Range("F8").Value = 51
$ws4.Range("F9").Value = 378
$ws4.Range("F10").Value = 7
